$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# "Scénario 2 Lecture d'une personne existante" gains a bold ".a" right
# after the bold "2", turning the heading into "Scénario 2.a Lecture ...".
$rng = $d.Content
$found = $rng.Find.Execute("Scénario 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Collapse to the point right after "2" (end of the found match) and
    # insert the new bold ".a" text there, before the following space run.
    $rng.Collapse(0)
    $rng.InsertAfter(".a")
    $rng.Font.Bold = 1
}

# --- Change 2 --------------------------------------------------------
# "Scénario 3 Lecture d'une personne inexistante" becomes
# "Scénario 2.b Lecture d'une personne inexistante" - the bold "3" run's
# text is replaced with "2.b".
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Scénario 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    # Narrow the match down to just the trailing "3" character and swap
    # its text for "2.b", preserving the bold formatting of that run.
    $threeRng = $d.Range($rng2.End - 1, $rng2.End)
    $threeRng.Text = "2.b"
}
